# Reporting Financial Calculations - update financial calculations with latest work progress
#
# This updates the "Sheet1" tab (the quarterly NeCTAR financial calculation sheet):
#  - Relabels the quarterly reporting period headers (row 1) to reflect the
#    latest reporting schedule.
#  - Updates the percentage-complete inputs (column L) for the 3 financial
#    milestones that have progressed, which drives the dependent
#    forecast/actual formulas (columns M, N, P, Q) and the summary totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Activate()

# Update quarter period labels in the header row.
$ws.Range("M1").Value = "March-June"
$ws.Range("P1").Value = "July-Sep"
$ws.Range("S1").Value = "Sep - Dec"

# Update progress percentages driving the cash-flow forecast calculations.
$ws.Range("L11").Value = 0.75
$ws.Range("L13").Value = 1
$ws.Range("L15").Value = 0.25

# Reflect latest view state on this sheet.
$ws.Range("S19").Select()
$ws.Application.ActiveWindow.Zoom = 88
